# Refresh the "想去人数" (interest count) figures in the 展览 (Exhibitions)
# and 全部类型 (All types) sheets to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 157
$ws1.Range("F5").Value = 1318
$ws1.Range("F6").Value = 18327
$ws1.Range("F7").Value = 375
$ws1.Range("F8").Value = 265
$ws1.Range("F10").Value = 6895
$ws1.Range("F14").Value = 117
$ws1.Range("F15").Value = 70
$ws1.Range("F19").Value = 254
$ws1.Range("F21").Value = 659
$ws1.Range("F22").Value = 37
$ws1.Range("F24").Value = 36
$ws1.Range("F25").Value = 279
$ws1.Range("F26").Value = 999
$ws1.Range("F27").Value = 130
$ws1.Range("F28").Value = 5174
$ws1.Range("F30").Value = 48
$ws1.Range("F32").Value = 77
$ws1.Range("F33").Value = 12130
$ws1.Range("F34").Value = 1289
$ws1.Range("F35").Value = 43
$ws1.Range("F36").Value = 211
$ws1.Range("F37").Value = 291
$ws1.Range("F38").Value = 3930

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 157
$ws4.Range("F5").Value = 1318
$ws4.Range("F6").Value = 18327
$ws4.Range("F7").Value = 375
$ws4.Range("F8").Value = 265
$ws4.Range("F10").Value = 6895
$ws4.Range("F14").Value = 117
$ws4.Range("F15").Value = 70
$ws4.Range("F19").Value = 254
$ws4.Range("F21").Value = 659
$ws4.Range("F22").Value = 37
$ws4.Range("F24").Value = 36
$ws4.Range("F25").Value = 279
$ws4.Range("F26").Value = 999
$ws4.Range("F27").Value = 130
$ws4.Range("F28").Value = 5174
$ws4.Range("F32").Value = 48
$ws4.Range("F34").Value = 77
$ws4.Range("F35").Value = 12130
$ws4.Range("F36").Value = 1289
$ws4.Range("F37").Value = 43
$ws4.Range("F38").Value = 211
$ws4.Range("F39").Value = 291
$ws4.Range("F40").Value = 3930
